# publikasi_template.xlsx edit
#
# Summary of the change (per the commit message "delete prodi using prodi
# from data mahasiswa/dosen instead" and the accompanying diff):
#   - every "Prodi*" column (Prodi, Prodi1..Prodi6, Prodi Lain) is removed
#     from the PUBLIKASI header row, since that info now comes from the
#     mahasiswa/dosen master data instead of being typed into this sheet;
#   - the trailing "... Lain" author block (Nama Penulis Lain / Status Lain /
#     Afiliasi Lain) is renumbered to a 7th numbered author slot
#     (Nama Penulis7 / Status7 / Afiliasi7), consistent with Nama Penulis1..6;
#   - everything to the right of each removed column shifts left, which is
#     why the sheet's dimension shrinks from A1:AL1 to A1:AD1 and every
#     column width shifts down by one column index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 8 "Prodi" columns (G, K, O, S, W, AA, AE, AI in the original
# 38-column layout). Removing them from right to left keeps each index
# valid for the next delete.
$prodiColumns = @(35, 31, 27, 23, 19, 15, 11, 7)
foreach ($colIndex in $prodiColumns) {
    $ws.Columns.Item($colIndex).Delete()
}

# The old "... Lain" (other/extra author) triple is now the 7th numbered
# author triple. After the deletions above these three headers land in
# columns AA, AB and AC.
$ws.Range("AA1").Value = "Nama Penulis7"
$ws.Range("AB1").Value = "Status7"
$ws.Range("AC1").Value = "Afiliasi7"

# Best-effort: restore the view's scroll position/selection (the sheet used
# to be scrolled so column AD was leftmost with AI4 selected; now it is
# scrolled to Q1 with AD2 selected).
$ws.Range("AD2").Select()
$excel.ActiveWindow.ScrollColumn = 17
